$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FoSYCRpUNL")

$newItems = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 19
foreach ($name in $newItems) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 0.005
    $row++
}

$ws.Range("A25").Select()

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
